$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 1) Header row: ARDUINO UNO / ARDUINO UNO 2 / ARDUINO MEGA -> ESCLAVO *
#    (single-run cells, so a direct Range.Text assignment is precise and
#    stays confined to the target cell)
# ---------------------------------------------------------------------------
$t.Cell(1, 2).Range.Text = "ESCLAVO 1"
$t.Cell(1, 3).Range.Text = "ESLCAVO 2"
$t.Cell(1, 4).Range.Text = "ESCLAVO 3"

# ---------------------------------------------------------------------------
# 2) "Entrada digital" row: the two middle "1" cells become "-"
# ---------------------------------------------------------------------------
$t.Cell(2, 3).Range.Text = "-"
$t.Cell(2, 4).Range.Text = "-"
$t.Cell(2, 5).Range.Text = "Se usará pulsador"

# ---------------------------------------------------------------------------
# 3) "Entrada análoga" row, description cell: rewrite the sentence. This
#    cell is made of several runs (it used to have proofErr spell-check
#    markers around "arduino"), so a plain Range.Text assignment would only
#    touch the first run and leave the rest behind. A document-level
#    Find/Replace operates on the rendered text across run boundaries, and
#    the sentence is unique in the document, so this is safe and replaces
#    every run that made up the old sentence with one clean run.
# ---------------------------------------------------------------------------
$oldText = "Para el arduino uno se usará un potenciómetro, para el mega un sensor de gas."
$newText = "Para el esclavo dos  se usará un potenciómetro, para el esclavo 3 un sensor de gas."
$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newText, 2)

# Move the _GoBack bookmark so it now sits right after "...esclavo 3" (it
# used to sit near the very end of the document, right before the closing
# picture).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$cell = $t.Cell(3, 5)
$findRange = $cell.Range
$null = $findRange.Find.Execute("para el esclavo 3", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$bookmarkPos = $findRange.End
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$null = $d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------------
# 4) "Salida digital" row: the ESLCAVO 2 column's "1" becomes "-"
# ---------------------------------------------------------------------------
$t.Cell(4, 3).Range.Text = "-"
